# Updated cryptos list (Price/Volume columns) to match latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell, forcing it to be stored as text even
# when it looks numeric (e.g. "93.06" or "30.009.17"), while leaving the
# cells formatting exactly as it was before (no left-over NumberFormat).
function Set-TextValue($cell, $text) {
    if ($text -match "^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$") {
        $cell.NumberFormat = "@"
        $cell.Value = $text
        $cell.ClearFormats()
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Range("D2") "30.009.17"
$ws.Range("E2").Value = "  -0.07%  "

Set-TextValue $ws.Range("D3") "1.910.09"
$ws.Range("E3").Value = "  +0.21%  "

Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("E5").Value = "  +3.94%  "

Set-TextValue $ws.Range("D6") "241.83"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("E7").Value = "  +0.07%  "

Set-TextValue $ws.Range("D8") "0.3159"
$ws.Range("E8").Value = "  +2.51%  "

Set-TextValue $ws.Range("D9") "26.32"
$ws.Range("E9").Value = "  +2.89%  "

Set-TextValue $ws.Range("D10") "0.06907"
$ws.Range("E10").Value = "  +0.01%  "

Set-TextValue $ws.Range("D11") "0.08003"
$ws.Range("E11").Value = "  -0.22%  "

Set-TextValue $ws.Range("D12") "1.908.79"
$ws.Range("E12").Value = "  +0.22%  "

Set-TextValue $ws.Range("D13") "0.7423"
$ws.Range("E13").Value = "  -1.75%  "

Set-TextValue $ws.Range("D15") "93.06"
$ws.Range("E15").Value = "  +1.32%  "

Set-TextValue $ws.Range("D16") "30.006.30"
$ws.Range("E16").Value = "  -0.08%  "

Set-TextValue $ws.Range("D17") "13.97"
$ws.Range("E17").Value = "  -0.63%  "

Set-TextValue $ws.Range("D18") "5.865"
$ws.Range("E18").Value = "  -5.26%  "

Set-TextValue $ws.Range("D19") "245.63"
$ws.Range("E19").Value = "  +3.29%  "

Set-TextValue $ws.Range("D20") "0.000007735"
$ws.Range("E20").Value = "  -0.26%  "

$ws.Range("E21").Value = "  +0.05%  "

Set-TextValue $ws.Range("D22") "2.152.49"
$ws.Range("E22").Value = "  +0.51%  "

$ws.Range("E23").Value = "  +0.14%  "

Set-TextValue $ws.Range("D24") "6.834"
$ws.Range("E24").Value = "  -3.00%  "

Set-TextValue $ws.Range("D25") "168.08"
$ws.Range("E25").Value = "  +0.87%  "

Set-TextValue $ws.Range("D26") "9.229"
$ws.Range("E26").Value = "  -0.99%  "

Set-TextValue $ws.Range("D27") "0.1388"
$ws.Range("E27").Value = "  +7.09%  "

Set-TextValue $ws.Range("D28") "18.94"
$ws.Range("E28").Value = "  +0.49%  "

Set-TextValue $ws.Range("D29") "2.032"
$ws.Range("E29").Value = "  -2.06%  "

$ws.Range("E30").Value = "  +1.53%  "

Set-TextValue $ws.Range("D31") "1.513"
$ws.Range("E31").Value = "  -0.93%  "

$ws.Range("E32").Value = "  +0.05%  "

Set-TextValue $ws.Range("D33") "4.082"
$ws.Range("E33").Value = "  +0.85%  "

Set-TextValue $ws.Range("D34") "0.05509"
$ws.Range("E34").Value = "  +1.97%  "

Set-TextValue $ws.Range("D35") "1.258"
$ws.Range("E35").Value = "  -2.41%  "

Set-TextValue $ws.Range("D36") "0.7317"
$ws.Range("E36").Value = "  -0.89%  "

Set-TextValue $ws.Range("D37") "2.721"
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("E38").Value = "  -1.14%  "

Set-TextValue $ws.Range("D39") "2.780"
$ws.Range("E39").Value = "  +0.55%  "

Set-TextValue $ws.Range("D40") "6.124"
$ws.Range("E40").Value = "  -2.07%  "

Set-TextValue $ws.Range("D41") "0.4411"
$ws.Range("E41").Value = "  -1.03%  "

Set-TextValue $ws.Range("D42") "72.27"
$ws.Range("E42").Value = "  -0.82%  "

$ws.Range("E43").Value = "  +0.05%  "

Set-TextValue $ws.Range("D44") "0.8373"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("E45").Value = "  -3.89%  "

$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D46") "7.546"
$ws.Range("E46").Value = "  -2.13%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D47") "100.31"
$ws.Range("E47").Value = "  -1.35%  "

Set-TextValue $ws.Range("D48") "986.62"
$ws.Range("E48").Value = "  +7.05%  "

Set-TextValue $ws.Range("D49") "2.057.90"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("E50").Value = "  -0.75%  "

$ws.Range("E51").Value = "  -0.60%  "
